$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 3 new blank rows right after row 79 (rows 80-82 become new/blank),
# pushing the existing content at rows 80+ down to rows 83+.
$ws.Rows("80:82").Insert()

# Row 79: record the new admin-capability entry (hours + description)
$ws.Cells.Item(79, 2).Value = 4.5
$ws.Range("D79").Value = "Added capability for an admin to unflag a story or block a user and hide their story.  "
$ws.Range("D79").Style = "Normal"

# Row 82: the "Sun am" note moves here (was row 80 before the insert)
$ws.Range("D82").Value = "Sun am - 3 - allows admin to unflag a story or block a user and hide their stories."

# Row 85: "For next deploy I WILL need new database!" (highlighted note, was row 83 before insert)
$ws.Range("D85").Value = "For next deploy I WILL need new database!"

# Row 86 gets a custom row height (30) and no longer carries the old "Need blocked user..." note
$ws.Rows(86).RowHeight = 30
$ws.Range("D86").ClearContents()

# Row 87: replace "Need blocked user to actually not be able to write" with the new follow-up note
$ws.Range("D87").Value = "For db - add read write ability to admin"

# Row 89 (was row 85 pre-insert): the "NB - the topic search..." note
$ws.Range("D89").Value = "NB - the topic search results really aren't very useful without a profile link - I should do it if I can manage it! Hahahahahaha"

# Row 90: "Next biggest issue is validation..." note
$ws.Range("D90").Value = "Next biggest issue is validation, I think, and of course interface."

# Update the view so the selection/scroll position matches the post-edit state
$ws.Application.ActiveWindow.ScrollRow = 72
$ws.Range("D79").Select()
